# Update "想去人数" (interested-people count) figures in column F
# for the "展览" and "全部类型" worksheets, matching the refreshed
# data snapshot recorded in the commit.

$wb = $excel.ActiveWorkbook

# row (by column A's "A<row>") -> new value for column F
$updates = @{
    2  = 294
    4  = 10239
    6  = 933
    7  = 1273
    8  = 6629
    10 = 433
    13 = 3151
    15 = 307
    16 = 629
    18 = 465
    20 = 52
    21 = 1595
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
